# Auto-generated edit script applying the 'Arribos 141 actualizados - 19' update
# across the TODOS, 215 and COMBINADAS sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('TODOS')
$ws.Cells.Item(2, 1).Value = '15:51'
$ws.Cells.Item(2, 2).Value = '11_ETCHEVERRY'
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(3, 1).Value = '15:51'
$ws.Cells.Item(3, 2).Value = '16_SANTA ANA'
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(4, 1).Value = '15:55'
$ws.Cells.Item(4, 2).Value = '17_ROMERO'
$ws.Cells.Item(4, 3).Value = 4
$ws.Cells.Item(4, 4).Value = '📅'
$ws.Cells.Item(5, 2).Value = '27_EL RETIRO'
$ws.Cells.Item(5, 3).Value = 5
$ws.Cells.Item(5, 4).Value = '🚌'
$ws.Cells.Item(6, 3).Value = 10
$ws.Cells.Item(7, 3).Value = 11
$ws.Cells.Item(8, 3).Value = 13
$ws.Cells.Item(9, 3).Value = 17
$ws.Cells.Item(10, 3).Value = 22
$ws.Cells.Item(11, 3).Value = 24
$ws.Cells.Item(12, 1).Value = '16:19'
$ws.Cells.Item(12, 3).Value = 28
$ws.Cells.Item(13, 3).Value = 30
$ws.Cells.Item(14, 3).Value = 38
$ws.Cells.Item(15, 3).Value = 39
$ws.Cells.Item(16, 3).Value = 45
$ws.Cells.Item(17, 3).Value = 51
$ws.Cells.Item(18, 3).Value = 52
$ws.Cells.Item(19, 3).Value = 57
$ws.Cells.Item(20, 1).Value = '16:52'
$ws.Cells.Item(20, 2).Value = '215B_LP-P MOR-40 Y 115'
$ws.Cells.Item(21, 1).Value = '16:56'
$ws.Cells.Item(21, 2).Value = '17_179 Y 38'
$ws.Cells.Item(21, 3).Value = 65
$ws.Cells.Item(22, 1).Value = '17:04'
$ws.Cells.Item(22, 2).Value = '215A_EL PATO'
$ws.Cells.Item(22, 3).Value = 73
$ws.Cells.Item(23, 1).Value = '17:07'
$ws.Cells.Item(23, 2).Value = '23_HERNANDEZ'
$ws.Cells.Item(23, 3).Value = 76
$ws.Cells.Item(24, 1).Value = '17:14'
$ws.Cells.Item(24, 2).Value = '215A_LA PLATA'
$ws.Cells.Item(24, 3).Value = 83
$ws.Cells.Item(25, 1).Value = '17:21'
$ws.Cells.Item(25, 2).Value = '26_HERNANDEZ'
$ws.Cells.Item(25, 3).Value = 90
$ws.Cells.Item(26, 1).Value = '17:24'
$ws.Cells.Item(26, 2).Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Cells.Item(26, 3).Value = 93
$ws.Cells.Item(27, 1).Value = '17:28'
$ws.Cells.Item(27, 2).Value = '14_ABASTO'
$ws.Cells.Item(27, 3).Value = 97
$ws.Cells.Item(28, 3).Value = 105
$ws.Cells.Item(29, 3).Value = 107
$ws.Cells.Item(30, 3).Value = 109
$ws.Cells.Item(31, 1).Value = '17:50'
$ws.Cells.Item(31, 2).Value = '16_P MOR-167 Y 521'
$ws.Cells.Item(31, 3).Value = 119
$ws.Cells.Item(31, 4).Value = '📅'

$ws = $wb.Worksheets.Item('215')
$ws.Cells.Item(2, 3).Value = 22
$ws.Cells.Item(3, 1).Value = '16:19'
$ws.Cells.Item(3, 3).Value = 28
$ws.Cells.Item(4, 1).Value = '16:52'
$ws.Cells.Item(4, 3).Value = 61
$ws.Cells.Item(5, 3).Value = 73
$ws.Cells.Item(6, 3).Value = 83
$ws.Cells.Item(7, 3).Value = 109

$ws = $wb.Worksheets.Item('COMBINADAS')
$ws.Cells.Item(2, 1).Value = '15:51'
$ws.Cells.Item(2, 2).Value = '11_ETCHEVERRY'
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(3, 1).Value = '15:51'
$ws.Cells.Item(3, 2).Value = '16_SANTA ANA'
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(4, 1).Value = '15:55'
$ws.Cells.Item(4, 2).Value = '17_ROMERO'
$ws.Cells.Item(4, 3).Value = 4
$ws.Cells.Item(4, 4).Value = '📅'
$ws.Cells.Item(5, 2).Value = '27_EL RETIRO'
$ws.Cells.Item(5, 3).Value = 5
$ws.Cells.Item(5, 4).Value = '🚌'
$ws.Cells.Item(6, 3).Value = 10
$ws.Cells.Item(7, 3).Value = 11
$ws.Cells.Item(8, 3).Value = 13
$ws.Cells.Item(9, 3).Value = 17
$ws.Cells.Item(10, 3).Value = 22
$ws.Cells.Item(11, 3).Value = 24
$ws.Cells.Item(12, 1).Value = '16:19'
$ws.Cells.Item(12, 3).Value = 28
$ws.Cells.Item(13, 3).Value = 30
$ws.Cells.Item(14, 3).Value = 38
$ws.Cells.Item(15, 3).Value = 39
$ws.Cells.Item(16, 3).Value = 45
$ws.Cells.Item(17, 3).Value = 51
$ws.Cells.Item(18, 3).Value = 52
$ws.Cells.Item(19, 3).Value = 57
$ws.Cells.Item(20, 1).Value = '16:52'
$ws.Cells.Item(20, 2).Value = '215B_LP-P MOR-40 Y 115'
$ws.Cells.Item(21, 1).Value = '16:56'
$ws.Cells.Item(21, 2).Value = '17_179 Y 38'
$ws.Cells.Item(21, 3).Value = 65
$ws.Cells.Item(22, 1).Value = '17:04'
$ws.Cells.Item(22, 2).Value = '215A_EL PATO'
$ws.Cells.Item(22, 3).Value = 73
$ws.Cells.Item(23, 1).Value = '17:07'
$ws.Cells.Item(23, 2).Value = '23_HERNANDEZ'
$ws.Cells.Item(23, 3).Value = 76
$ws.Cells.Item(24, 1).Value = '17:14'
$ws.Cells.Item(24, 2).Value = '215A_LA PLATA'
$ws.Cells.Item(24, 3).Value = 83
$ws.Cells.Item(25, 1).Value = '17:21'
$ws.Cells.Item(25, 2).Value = '26_HERNANDEZ'
$ws.Cells.Item(25, 3).Value = 90
$ws.Cells.Item(26, 1).Value = '17:24'
$ws.Cells.Item(26, 2).Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Cells.Item(26, 3).Value = 93
$ws.Cells.Item(27, 1).Value = '17:28'
$ws.Cells.Item(27, 2).Value = '14_ABASTO'
$ws.Cells.Item(27, 3).Value = 97
$ws.Cells.Item(28, 3).Value = 105
$ws.Cells.Item(29, 3).Value = 107
$ws.Cells.Item(30, 3).Value = 109
$ws.Cells.Item(31, 1).Value = '17:50'
$ws.Cells.Item(31, 2).Value = '16_P MOR-167 Y 521'
$ws.Cells.Item(31, 3).Value = 119
$ws.Cells.Item(31, 4).Value = '📅'

